$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 480) holds a "Förändrad" (changed) date serial.
# All of these values move from 45178 to 45179 (one day later).
$ws.Range("C2:C480").Value = 45179
